$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Item "balance" (الرصيد الحالي, column H) text updates -----------------
# Row 14: CONTAFEVER N 200MG/5ML SUSP. 120ML  "20:0" -> "19:0"
$ws.Range("H14").Value = "19:0"

# Row 18: DIASMECT 20% SUSP. 60ML  "16:0" -> "15:0"
$ws.Range("H18").Value = "15:0"

# Row 20: DOLIPRANE 1 GM 15 TABS.  "6:2" -> "6:1"
$ws.Range("H20").Value = "6:1"

# Row 22: GASTROMOTIL 1MG/ML ORAL SUSP. 200ML  "2:0" -> "1:0"
$ws.Range("H22").Value = "1:0"

# Row 26: MOTILIUM 10MG 40 F.C.TAB.  "0:1" -> "0:0"
$ws.Range("H26").Value = "0:0"

# --- Row 32 (STREPTOQUIN 20 TABLETS) updates --------------------------------
# P32 (sale price) "23.0000" -> "0.0000" -- keep stored as text, matching the
# source data (which stores these numeric-looking values as literal text),
# and keep the existing "0.00" number-format style on the cell.
$p32Format = $ws.Range("P32").NumberFormat
$ws.Range("P32").NumberFormat = "@"
$ws.Range("P32").Value = "0.0000"
$ws.Range("P32").NumberFormat = $p32Format

# Q32 (number of transactions) "0:1" -> "0:0"
$ws.Range("Q32").Value = "0:0"

# --- Grand total (row 38, column P) -----------------------------------------
# Literal total of column P; STREPTOQUIN's sale price dropped by 23.00
# (23.0000 -> 0.0000), so the total drops by the same amount.
$ws.Range("P38").Value = 1741.1300000000001
